$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 58: "Cas na vyplneni - zmeneno na sekundy" ---
$ws.Range("B58").Value = "Cas na vyplneni - zmeneno na sekundy"
$ws.Range("C58").Value = 1.5
$ws.Range("D57").Copy()
$ws.Range("D58").PasteSpecial(-4122)
$ws.Range("D58").Value = 40983

# --- Row 59: "Maximalni pocet vyplneni, ..." ---
$ws.Range("B59").Value = "Maximalni pocet vyplneni, mazani otazky, mazani alternativni otazky, skupiny pri editaci formu, styly, michani choice odpovedi"
$ws.Range("C59").Value = 6
$ws.Range("D57").Copy()
$ws.Range("D59").PasteSpecial(-4122)
$ws.Range("D59").Value = 40986

# --- Row 60: "Google services studium, ..." ---
$ws.Range("B60").Value = "Google services studium, Google Groups csv import, vlastni obsluha vyjimek"
$ws.Range("C60").Value = 5
$ws.Range("D57").Copy()
$ws.Range("D60").PasteSpecial(-4122)
$ws.Range("D60").Value = 40987

# Update the selection to match post-edit cursor position
$ws.Range("D61").Select()
